# Cross-walk update after re-running the fuzzy matcher.
# The fuzzy match for this facility changed from MAYFAIR ELEMENTARY SCHOOL / Kenderton
# to Mary McLeod Bethune School / Mary M. Bethune, so refresh the matched row (row 2)
# with the new facility/school names, addresses and match probability.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "Mary McLeod Bethune School"
$ws.Range("B2").Value = "Mary M. Bethune School"
$ws.Range("C2").Value = "3301 Old York Road"
$ws.Range("D2").Value = "3301 OLD YORK RD"
$ws.Range("E2").Value = 0.89

# The crosswalk sheet no longer needs the wide, bestFit columns C:F (their data stays,
# only the explicit widths are dropped back toward the sheet default) while A and B get
# new explicit widths sized for the longer re-matched names.
$ws.Columns.Item(1).ColumnWidth = 23
$ws.Columns.Item(2).ColumnWidth = 22.15
$ws.Columns.Item(3).ColumnWidth = 8
$ws.Columns.Item(4).ColumnWidth = 8
$ws.Columns.Item(5).ColumnWidth = 8
$ws.Columns.Item(6).ColumnWidth = 8
